$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Option A:" -> "Option A"  (drop trailing colon)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Option A:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Option A", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. "Option B: " -> "Option B " (drop colon, keep trailing space)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Option B: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Option B ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3. "Distance from you" -> "Distance from You", ending up split across
#    three runs ("Distance from " / "Y" / "ou") like the target revision.
#    A plain Find/Replace collapses same-formatted runs in the edited
#    paragraph into one run, so after the text fix we nudge Bold on/off
#    (net no-op) over each sub-range to force the engine to re-split the
#    run boundaries at the desired offsets.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Distance from you", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Distance from You", 2) | Out-Null

$found = $d.Content
$found.Find.Execute("Distance from You", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$start = $found.Start
foreach ($seg in @(@(0,14), @(14,15), @(15,17))) {
    $r = $d.Range($start + $seg[0], $start + $seg[1])
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------
# 4. "cost to you" -> "Cost to You", ending up split across four runs
#    ("C" / "ost to " / "Y" / "ou") like the target revision, while the
#    preceding "One-time " run is left completely untouched.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("cost to you", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Cost to You", 2) | Out-Null

$found2 = $d.Content
$found2.Find.Execute("Cost to You", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$start2 = $found2.Start
foreach ($seg in @(@(0,1), @(1,8), @(8,9), @(9,11))) {
    $r = $d.Range($start2 + $seg[0], $start2 + $seg[1])
    $r.Bold = 1
    $r.Bold = 0
}

# ---------------------------------------------------------------------
# 5. Footer: text is unchanged, but the eight runs collapse into a
#    single run sharing the formatting common to all of them.
# ---------------------------------------------------------------------
$footerRange = $d.Sections(1).Footers(1).Range
$footerRange.Find.Execute(
    "If Option A was the ONLY option besides No Project, which would you choose?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "If Option A was the ONLY option besides No Project, which would you choose?", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Header: remove the "Question Number «MERGEFIELD card»" field text
#    and replace with plain text "Two options, one choice".
# ---------------------------------------------------------------------
$headerRange = $d.Sections(1).Headers(1).Range
$headerRange.Fields(1).Delete()
$headerRange.Find.Execute("Question Number ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "Two options, one choice", 2) | Out-Null
